$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 359; all existing rows 359..427 shift down to 360..428
$ws.Rows.Item(359).Insert()

# Populate the newly inserted row 359 with the new weekly record
$ws.Cells.Item(359,1).Value  = 9
$ws.Cells.Item(359,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(359,3).Value  = "Metropolitana"
$ws.Cells.Item(359,4).Value  = 44694
$ws.Cells.Item(359,5).Value  = 13
$ws.Cells.Item(359,6).Value  = 100112012
$ws.Cells.Item(359,7).Value  = "Espinaca"
$ws.Cells.Item(359,8).Value  = "Sin especificar"
$ws.Cells.Item(359,9).Value  = "Primera"
$ws.Cells.Item(359,10).Value = 160
$ws.Cells.Item(359,11).Value = 6000
$ws.Cells.Item(359,12).Value = 7000
$ws.Cells.Item(359,13).Value = 6500
$ws.Cells.Item(359,14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(359,15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(359,16).Value = 650
$ws.Cells.Item(359,17).Value = 10
$ws.Cells.Item(359,18).Value = "Hortaliza"
